# "added a reference to the administrative simplification slides"
#
# 1. Slide 5 ("Example - GIC Tiering"): the second paragraph of the body
#    placeholder was re-typed so its text is one contiguous run instead of
#    two runs that happen to concatenate to the same string.
# 2. Slide 6 ("References"): a new reference (FCHP GIC tiering page) was
#    added as its own paragraph/hyperlink between the mass.gov reference and
#    the betterhealthconnector.com reference, and the body placeholder was
#    switched to "Shrink text on overflow" (normAutofit) now that there is
#    one more line of text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5: merge the two runs of paragraph 2 back into a single run
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$para2 = $tr5.Paragraphs(2, 1)
# First retype with a throwaway value so the run-merge actually happens
# (setting identical text is treated as a no-op by the host).
$para2.Text = "__tmp__"
$tr5b = $sh5.TextFrame.TextRange
$para2b = $tr5b.Paragraphs(2, 1)
$para2b.Text = "This ranking process uses detailed claims data.  Historically, this required insurance companies to submit data directly to the GIC."

# ---------------------------------------------------------------------
# Slide 6: add the FCHP reference paragraph + hyperlink
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)

# Shrink text on overflow, since the placeholder now holds one more line.
$sh6.TextFrame.AutoSize = 2

$tr6 = $sh6.TextFrame.TextRange

# Free up the relationship id currently used by the betterhealthconnector.com
# link so that the newly-added FCHP link claims it instead (matches the
# target file, where the new link is rId4 and betterhealthconnector.com
# becomes rId5).
$paraBhc = $tr6.Paragraphs(3, 1)
$paraBhc.ActionSettings(1).Hyperlink.Address = ""

# Insert the new paragraph right after the mass.gov reference (paragraph 2).
$tr6b = $sh6.TextFrame.TextRange
$paraMassGov = $tr6b.Paragraphs(2, 1)
$fchpUrl = "http://www.fchp.org/providers/resources/gic-tiering.aspx"
[void]$paraMassGov.InsertAfter("`r" + $fchpUrl)

# Re-fetch the new paragraph (now paragraph 3) and wire up its hyperlink,
# split into the same run layout as the other references ("http", "://",
# rest of the address).
$tr6c = $sh6.TextFrame.TextRange
$newPara = $tr6c.Paragraphs(3, 1)
$newStart = $newPara.Start

$run1 = $tr6c.Characters($newStart, 4)
$run1.ActionSettings(1).Hyperlink.Address = $fchpUrl

$tr6d = $sh6.TextFrame.TextRange
$run2 = $tr6d.Characters($newStart + 4, 3)
$run2.ActionSettings(1).Hyperlink.Address = $fchpUrl

$tr6e = $sh6.TextFrame.TextRange
$run3 = $tr6e.Characters($newStart + 7, 49)
$run3.ActionSettings(1).Hyperlink.Address = $fchpUrl

# Restore the betterhealthconnector.com hyperlink (now paragraph 4); this
# claims the next free relationship id (rId5).
$tr6f = $sh6.TextFrame.TextRange
$paraBhc2 = $tr6f.Paragraphs(4, 1)
$paraBhc2.ActionSettings(1).Hyperlink.Address = "https://betterhealthconnector.com/wp-content/uploads/board_meetings/2015/2015-07-09/Board-Memo-Risk-Adjustment-Update-070615.pdf"
